# Apply the "productItems" commit: extend product/products sheets with new
# rows, add header + wrap-text styling, add a new "productItems" sheet, and
# append the new shared strings (Airpods, PS4, Mivi ...).

$wb = $excel.ActiveWorkbook

$wsProduct = $wb.Worksheets.Item("product")
$wsProducts = $wb.Worksheets.Item("products")

# Snapshot "product" (still pristine at this point) into a new sheet that
# will become "productItems" - this keeps the same xml namespaces/column
# defaults the original authored sheets use instead of a bare-bones blank
# sheet. Its content gets wiped and rebuilt below.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsProduct.Copy($null, $lastSheet)
$wsItems = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsItems.Name = "productItems"
$wsItems.Cells.ClearContents()

# ---------------------------------------------------------------------
# 1) Sheet "product" (sheet1) - fill in the new rows first so the shared
#    strings get created in the same order the target file uses them
#    (Airpods=4, PS4=5, Mivi...=6).
# ---------------------------------------------------------------------
$wsProduct.Range("A3").Value = "Samsung Galaxy A1"
$wsProduct.Range("A4").Value = "Airpods"
$wsProduct.Range("A5").Value = "PS4"
$wsProduct.Range("A6").Value = "Mivi Play 5 Watt Truly Wireless Bluetooth Portable Speaker (Black)"

# Header style: bold font + yellow fill.
$wsProduct.Range("A1").Font.Bold = $true
$wsProduct.Range("A1").Interior.Color = 65535

# Wrap-text style used for the long description cell.
$wsProduct.Range("A6").Font.Color = 1118479
$wsProduct.Range("A6").WrapText = $true
$wsProduct.Range("A6").VerticalAlignment = -4108
$wsProduct.Rows.Item(6).RowHeight = 30

$wsProduct.Columns.Item(1).ColumnWidth = 34.333333333333336

$wsProduct.PageSetup.PaperSize = 9
$wsProduct.PageSetup.Orientation = 1

$wsProduct.Range("A6").Select() | Out-Null

# ---------------------------------------------------------------------
# 2) Sheet "products" (sheet2) - reuse the styles created above via
#    copy/paste-special so no orphan style entries get created.
# ---------------------------------------------------------------------
$wsProducts.Range("A3").Value = "Samsung Galaxy A1"
$wsProducts.Range("B3").Value = "Airpods"
$wsProducts.Range("A4").Value = "PS4"
$wsProducts.Range("B4").Value = "Samsung Galaxy A1"
$wsProducts.Range("A5").Value = "Mivi Play 5 Watt Truly Wireless Bluetooth Portable Speaker (Black)"
$wsProducts.Range("B5").Value = "Airpods"
$wsProducts.Range("A6").Value = "Mivi Play 5 Watt Truly Wireless Bluetooth Portable Speaker (Black)"
$wsProducts.Range("B6").Value = "Hair Dryer"

$wsProduct.Range("A1").Copy()
$wsProducts.Range("A1").PasteSpecial(-4122)
$wsProducts.Range("B1").PasteSpecial(-4122)

$wsProduct.Range("A6").Copy()
$wsProducts.Range("A5").PasteSpecial(-4122)
$wsProducts.Range("A6").PasteSpecial(-4122)

$wsProducts.Rows.Item(5).RowHeight = 60
$wsProducts.Rows.Item(6).RowHeight = 60

$wsProducts.Columns.Item(1).ColumnWidth = 17.333333333333332

$wsProducts.Range("A6").Select() | Out-Null

# ---------------------------------------------------------------------
# 3) New sheet "productItems" (sheet3, created above from the "product"
#    snapshot) - fill in its own values/styling.
# ---------------------------------------------------------------------
$wsItems.Range("A1").Value = "Hair Dryer"
$wsItems.Range("A2").Value = "Samsung Galaxy A1"
$wsItems.Range("A3").Value = "Airpods"
$wsItems.Range("A5").Value = "PS4"

$wsProduct.Range("A6").Copy()
$wsItems.Range("A4").PasteSpecial(-4122)

$wsItems.Columns.Item(1).ColumnWidth = 17.333333333333332

$wsItems.Range("A5").Select() | Out-Null

# Restore "products" as the active sheet (it was active before the edit;
# adding/copying sheets switches focus away from it).
$wsProducts.Activate() | Out-Null
